# Update Name of Algo
# Apply updated imputed values to the result_data_KNN sheet (columns A and B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"  = -20.515
    "A7"  = -21.286
    "B7"  = 6.312
    "B15" = 5.042
    "A16" = -21.965
    "B21" = 9.000000000000002
    "B22" = 7.034000000000001
    "B23" = 7.359999999999999
    "A28" = -21.918
    "A29" = -21.586
    "A32" = -21.768
    "B34" = 8.059000000000001
    "A40" = -20.481
    "B43" = 5.448
    "B45" = 5.672000000000001
    "B50" = 5.366
    "B51" = 6.464
    "A52" = -21.684
    "A57" = -22.275
    "A66" = -21.504
    "B66" = 5.523
    "B67" = 5.571
    "B79" = 5.571999999999999
    "B84" = 5.453000000000001
    "B92" = 5.495000000000001
    "B97" = 6.156
    "A100" = -21.927
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
